$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.890759587287903
$ws.Range("B1").Value = 2.626100063323975
$ws.Range("C1").Value = 2.861844301223755
$ws.Range("D1").Value = 3.013024568557739
$ws.Range("E1").Value = 0.9806374907493591
